$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Eric moved his data from drive F: to drive E:
$ws.Range("I3").Value = "E:\Eric\Local\07_Transcription\RawData"
$ws.Range("I4").Value = "E:\Eric\Local\07_Transcription\FISHAnalysisData"
$ws.Range("I5").Value = "E:\Eric\Dropbox\Lab\07_Transcription\LivemRNAData"
$ws.Range("I12").Value = "E:\Eric\GitHub\mRNADynamics"

# Update the view's frozen pane location and active selection cell
$ws.Activate()
$ws.Range("I12").Select()
